$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F1").Value = "Banco"
$ws.Range("F2").Value = "Bancolombia"
$ws.Range("F3").Value = "Itau"
$ws.Range("F4").Value = "Bancolombia"
$ws.Range("F5").Value = "Itau"
$ws.Range("F6").Value = "Davivienda"
